# TestArrayAggregations.xlsx - introduce new methods to the Function
# Adds a new column G (id 1/2/3 cycling) next to the existing data table,
# replaces the old row 18/19 SUBSTITUTE() formulas with a new non-array
# F18 aggregation, a new G18 array-constant SUM(), and a new F19 MATCH()
# array formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: cycling 1,2,3 values for rows 2-17 -----------------
for ($r = 2; $r -le 17; $r++) {
    $val = (($r - 2) % 3) + 1
    $ws.Cells.Item($r, 7).Value = $val
}

# --- Clear the old row 18 / row 19 formulas ----------------------------
$ws.Range("F18").ClearContents()
$ws.Range("F19").ClearContents()

# --- Row 18: F18 gets the style from F17 (matches the rest of column F)
$ws.Range("F17").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").FormulaArray = "=D18:D33*E18:E33"

# G18: array-constant SUM formula
$ws.Range("G18").FormulaArray = "=SUM({1,2,3;1,2,3;1,2,3})"

# --- Row 19: F19 array MATCH() formula ---------------------------------
$ws.Range("F19").FormulaArray = "=MATCH(13300,F2:F17,0)"

# --- Fix up the selection so it matches the new single-cell reference -
$ws.Range("F18").Select()

$excel.CutCopyMode = 0
